$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the two new columns (I = "I0", J = "IF"),
# matching the style already used by the existing header row (row 1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2..88: I<row> and J<row> numeric values.
$rowData = @(
    @(2, 5, 6),
    @(3, 5, 6),
    @(4, 6, 6),
    @(5, 7, 7),
    @(6, 6, 7),
    @(7, 6, 6),
    @(8, 6, 6),
    @(9, 6, 7),
    @(10, 6, 6),
    @(11, 6, 7),
    @(12, 6, 6),
    @(13, 7, 7),
    @(14, 5, 5),
    @(15, 6, 6),
    @(16, 8, 9),
    @(17, 7, 7),
    @(18, 7, 7),
    @(19, 7, 7),
    @(20, 7, 7),
    @(21, 8, 8),
    @(22, 5, 6),
    @(23, 6, 6),
    @(24, 7, 7),
    @(25, 7, 7),
    @(26, 6, 6),
    @(27, 6, 6),
    @(28, 7, 7),
    @(29, 7, 7),
    @(30, 6, 6),
    @(31, 8, 8),
    @(32, 6, 7),
    @(33, 6, 7),
    @(34, 6, 6),
    @(35, 8, 8),
    @(36, 5, 6),
    @(37, 10, 10),
    @(38, 8, 8),
    @(39, 7, 7),
    @(40, 5, 5),
    @(41, 9, 9),
    @(42, 7, 7),
    @(43, 7, 7),
    @(44, 7, 7),
    @(45, 8, 8),
    @(46, 6, 6),
    @(47, 8, 8),
    @(48, 5, 5),
    @(49, 6, 7),
    @(50, 9, 9),
    @(51, 5, 6),
    @(52, 4, 4),
    @(53, 7, 7),
    @(54, 7, 7),
    @(55, 6, 6),
    @(56, 7, 7),
    @(57, 8, 8),
    @(58, 6, 6),
    @(59, 1, 3),
    @(60, 1, 3),
    @(61, 4, 5),
    @(62, 6, 6),
    @(63, 7, 7),
    @(64, 3, 5),
    @(65, 6, 6),
    @(66, 4, 4),
    @(67, 5, 6),
    @(68, 7, 7),
    @(69, 9, 9),
    @(70, 1, 2),
    @(71, 3, 4),
    @(72, 8, 8),
    @(73, 3, 4),
    @(74, 8, 8),
    @(75, 3, 3),
    @(76, 3, 4),
    @(77, 5, 5),
    @(78, 7, 7),
    @(79, 8, 8),
    @(80, 7, 7),
    @(81, 9, 9),
    @(82, 7, 8),
    @(83, 9, 9),
    @(84, 8, 8),
    @(85, 6, 6),
    @(86, 7, 7),
    @(87, 4, 4),
    @(88, 4, 4)
)

foreach ($entry in $rowData) {
    $r = $entry[0]
    $iv = $entry[1]
    $jv = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iv
    $ws.Cells.Item($r, 10).Value = $jv
}
